# Add a new "VSC Trunking" setting row to the "Common" configuration sheet,
# right before the existing "Security" section, and document it with a
# comment (matching the author/description used for the other settings).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Common")

# Insert a new row at 60. This shifts the "Security" section header (and
# every row below it) down by one row, and Excel automatically re-points
# merged cells, data validations and existing cell comments to follow the
# shift.
$ws.Rows.Item(60).Insert()

# The freshly inserted row copies A60's look from the row above (a normal
# label cell), but B60 ends up with a slightly different blended style
# because the row that used to be 60 ("Security") had no B cell. Re-copy
# the standard "value cell" formatting from a neighboring data row so B60
# matches the other fields (e.g. B59, "OpenStack CA Certificate").
$ws.Range("B59").Copy()
$ws.Range("B60").PasteSpecial(-4122)

# Populate the new row with its label.
$ws.Range("A60").Value = "VSC Trunking"

# Document the new field with a comment, just like the other settings.
$ws.Range("A60").AddComment("This enables trunking between vsc control port and its underlay ports (openstack_external_port1_name, openstack_external_port2_name,openstack_external_port3_name) as defined under vscs.yml [default: False]")

# Give the new field the same true/false validation used by similar flags
# (e.g. XMPP TLS / OpenFlow TLS) just below it.
$ws.Range("B60").Validation.Add(3, 2, 1, "true,false")
$validation = $ws.Range("B60").Validation
$validation.ErrorTitle = "Invalid Entry"
$validation.ErrorMessage = "Your entry is not true or false, change anyway?"
$validation.InputTitle = "True or False Selection"
$validation.InputMessage = "Please select true or false"
